$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数") on sheets "展览" and "全部类型"
$updates = @{
    2  = 161
    3  = 422
    4  = 12309
    6  = 136
    10 = 190
    16 = 365
    17 = 3220
    19 = 935
    20 = 14
    22 = 29
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
